$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 0, $null),
    @(5, 13, $null),
    @(6, $null, "No answer"),
    @(7, 12, $null),
    @(8, 12, $null),
    @(9, 10, $null),
    @(10, 8, $null),
    @(11, $null, "Can't load image"),
    @(12, 16, $null),
    @(13, -1, "Error in the system , error (0)"),
    @(14, -1, "Error in the system , error (0)"),
    @(15, -1, "Error in the system , error (0)"),
    @(16, -1, "Error in the system , error (0)"),
    @(17, $null, "No answer"),
    @(18, -1, "Error in the system , error (0)"),
    @(19, -1, "Error in the system , error (0)"),
    @(20, -1, "Error in the system , error (0)"),
    @(21, 12, $null),
    @(22, 12, $null),
    @(23, 13, $null),
    @(24, 12, $null),
    @(25, 0, $null),
    @(26, 12, $null),
    @(27, 0, $null),
    @(28, 12, $null),
    @(29, 0, $null),
    @(30, 0, $null),
    @(31, 0, $null),
    @(32, 0, $null),
    @(33, 13, $null),
    @(34, 13, $null),
    @(35, 0, $null),
    @(36, 12, $null),
    @(37, 12, $null),
    @(38, 12, $null),
    @(39, 13, $null),
    @(40, 12, $null),
    @(41, $null, "No answer"),
    @(42, $null, "No answer"),
    @(43, $null, "No answer"),
    @(44, $null, "No answer"),
    @(45, $null, "Can't load image"),
    @(46, $null, "No answer"),
    @(47, 13, $null),
    @(48, 12, $null),
    @(49, 12, $null),
    @(50, 16, $null),
    @(51, 8, $null),
    @(52, 8, $null),
    @(53, 8, $null)
)

foreach ($row in $data) {
    $r = $row[0]
    $bVal = $row[1]
    $cVal = $row[2]
    if ($null -ne $bVal) {
        $ws.Cells.Item($r, 2).Value = $bVal
    }
    if ($null -ne $cVal) {
        $ws.Cells.Item($r, 3).Value = $cVal
    }
}

Write-Host "Applied $($data.Count) rows"
